$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append two new trading-log rows to the "intraday" sheet (rows 160/161).
#    Write column E before column B on each row so the new shared strings
#    land in the same append order as the target workbook.
# ---------------------------------------------------------------------------
$intraday = $wb.Worksheets.Item("intraday")

$intraday.Range("E160").Value = "I missed one big trade and againg taken two trade later both hit my SL so it's a good try anyway"
$intraday.Range("B160").Value = "first trade was ok but second trade is gambling so no gambling  if you trade too always follow your rule"

# ---------------------------------------------------------------------------
# 2) Add the new "Days" worksheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$days = $wb.Worksheets.Add($null, $lastSheet)
$days.Name = "Days"

$days.Range("A1").Value = "Days"
$days.Range("B1").Value = "Date"
$days.Range("C1").Value = "profit"
$days.Range("D1").Value = "Loss"

$days.Range("A2").Value = "Monday"
$days.Range("B2").Value = 45775
$days.Range("C2").Value = 0
$days.Range("D2").Value = 374

$days.Range("A3").Value = "Tuesday"
$days.Range("B3").Value = 45776
$days.Range("C3").Value = 0
$days.Range("D3").Value = 640

$days.Range("A4").Value = "Wednesday"
$days.Range("B4").Value = 45777
$days.Range("C4").Value = 0
$days.Range("D4").Value = 655

# Column widths for the new sheet (matches the authored workbook).
$days.Columns.Item(1).ColumnWidth = 11.17
$days.Columns.Item(2).ColumnWidth = 11.5
$days.Columns.Item(3).ColumnWidth = 11.5

# Match the date formatting already used elsewhere in the workbook by
# copying the number format from an existing date cell (reuses the same
# style record instead of minting a new one).
$intraday.Range("A159").Copy()
$days.Range("B2:B4").PasteSpecial(-4122)  # xlPasteFormats
$days.Range("B2").Value = 45775
$days.Range("B3").Value = 45776
$days.Range("B4").Value = 45777
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Finish the two new "intraday" rows (second mistake entry, then values).
# ---------------------------------------------------------------------------
$intraday.Range("E161").Value = "again I made two mistake first one was ok as per risk but second trade is not that satisfying control your psychology my marking was perfect in second trade but I entered lately as per chart comparison so try to enter at same point as market move 80 points down "
$intraday.Range("B161").Value = "only single trade no more two trade again don't break rule "

$intraday.Range("A159").Copy()
$intraday.Range("A160:A161").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$intraday.Range("A160").Value = 45776
$intraday.Range("C160").Value = 640
$intraday.Range("D160").Value = 2

$intraday.Range("A161").Value = 45777
$intraday.Range("C161").Value = 655
$intraday.Range("D161").Value = 2

# ---------------------------------------------------------------------------
# 4) Restore the selections/views to match the saved workbook state.
# ---------------------------------------------------------------------------
$retracement = $wb.Worksheets.Item("retracement_example")
$retracement.Select() | Out-Null
$retracement.Range("A13").Select() | Out-Null

$days.Select() | Out-Null
$days.Range("E4").Select() | Out-Null

$intraday.Select() | Out-Null
$intraday.Range("B161").Select() | Out-Null
